# summer 24 week 12 inputs
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$ws.Range("D2").Value = 11.27
$ws.Range("D3").Value = 10.23
$ws.Range("B4").Value = 8.73
$ws.Range("C4").Value = 9.77
$ws.Range("F5").Value = 10.2
$ws.Range("E6").Value = 9.800000000000001
$ws.Range("G6").Value = 10.35
$ws.Range("H6").Value = 10.53
$ws.Range("F7").Value = 9.65
$ws.Range("F8").Value = 9.470000000000001
$ws.Range("I8").Value = 8.880000000000001
$ws.Range("J8").Value = 11.4
$ws.Range("H9").Value = 11.12
$ws.Range("H10").Value = 8.6
